# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2 through 22
$kValues = @(1, 1, 1, 1, 3, 3, 1, 0, 2, 2, 2, 0, 0, 2, 2, 0, 3, 2, 1, 1, 1)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
